Write-Output "hi"
